# repull data, push all data, mean calculation
# Update the dSF column (F) for a set of rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = 2
    14 = 0
    15 = -1
    21 = 3
    22 = 0
    31 = -1
    33 = 3
    35 = -1
    44 = 1
    46 = 4
    49 = 1
    50 = 1
    54 = 0
    60 = 1
    64 = 6
    65 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
